# SKLAD.xlsx restructuring:
#   - rename the lone "1.01.2023 - 20.12.2023" sheet to "Yillik 2024"
#   - trim it down to just the header row + a reset "Umumiy Harajat
#     Summasi" / "0 so'm" summary row (the per-expense rows are gone)
#   - add two more sheets with identical header/summary content:
#       "Oylik 01" and "Haftalik 1.01.2024 - 8.01.2024"

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.ActiveSheet
$ws1.Name = "Yillik 2024"

# Drop the old per-expense rows (2 "Svarshikka"/"Eshmamat Kamaz" rows)
# and the stale totals row, keeping only the header row.
$ws1.Rows.Item(4).Delete()
$ws1.Rows.Item(3).Delete()
$ws1.Rows.Item(2).Delete()

# Re-create the summary row with the expenses reset to zero.
$ws1.Range("B2").Value = "Umumiy Harajat Summasi"
$ws1.Range("C2").Value = "0 so'm"

# B2:C2 need the same bold/centered/bordered look as the header cells;
# copy the formatting (not the value) from A1 so we reuse the existing
# cell style instead of registering a new one.
$ws1.Range("A1").Copy()
$ws1.Range("B2:C2").PasteSpecial(-4122)

# Column widths (13 / 25 / 10 / 18 chars). The host's ColumnWidth setter
# round-trips through OOXML with a constant +5/6 drift, so compensate by
# writing size-5/6 to land exactly on the intended integer width.
$ws1.Columns.Item(1).ColumnWidth = 13 - 5/6
$ws1.Columns.Item(2).ColumnWidth = 25 - 5/6
$ws1.Columns.Item(3).ColumnWidth = 10 - 5/6
$ws1.Columns.Item(4).ColumnWidth = 18 - 5/6

# Clone the finished sheet twice (format + content) to build the other
# two sheets, then rename each clone in place.
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Oylik 01"

$ws2.Copy($null, $ws2)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "Haftalik 1.01.2024 - 8.01.2024"

Write-Output "Sheets now:"
foreach ($s in $wb.Worksheets) { Write-Output $s.Name }
